$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Récupère la réponse du serveur avec l'id en paramètre sous format JSON"
$ws.Range("F5").Value = "La requète pourrait ne pas avoir le bon format, ne pas fonctionner ou pourrait retourner des informations erronées"
$ws.Range("D6").Value = "Récupère la couleur selectionnée par l'utilisateur"
$ws.Range("D7").Value = "Récupère la quantité selectionnée par l'utilisateur"
$ws.Range("D8").Value = "Récupère les informations lié à l'article et à sa personnalisation pour les stocker dans le locale storage et afficher une pop up confirmant l'ajout de l'article au panier , au clic sur le bouton ""ajouter au panier"""
$ws.Range("E8").Value = "Appeller la fonction avec différentes valeurs de tests (ex: ajouter différent articles ou ajouter différentes couleurs et quantités au même arrticle), et on observe la valeur retournée, avec un console.log(selectedArticle), puis vérifier que ces éléments soient bien stockés dans le locale storage via l'onglet Application du devtool"
$ws.Range("D9").Value = "Ferme la pop up lorsque l'utilisateur clique sur le bouton X"
$ws.Range("D10").Value = "Récupère les informations concernant les articles contenus dans locale storage"
$ws.Range("E10").Value = "Consulter le locale storage dans l'onglet application de l'outil devtool"
$ws.Range("D11").Value = "Une fois la page chargée, afficher la liste des articles contenus dans le locale storage (bloucle for ...of)`n2-calculer le prix total du panier `n3-afficher dynamiquement les informations liés aux articles ajoutés au panier"
$ws.Range("E11").Value = "Ajouter plusieurs articles différents au panier et vérifier qu'ils s'affichent bien et que les montants s'actualisent bien sur la page panier"
$ws.Range("F11").Value = "Tous les articles ajoutés au panier ne s'affichent pas, le prix total du panier ne s'affiche pas (si par exemple les données sont de type string et pas number) ou ne s'actualise pas (s'il n'y a pas de boucle permettant d'additionner le prix de chaque article du panier par exemple)"
$ws.Range("D12").Value = "Suprimme la ligne de l'article concerné sur la page`nweb"
$ws.Range("F12").Value = "La ligne pourrait ne pas être supprimée. Il faut également veiller à ce les autre lignes non concernées ne soient pas supprimées"
$ws.Range("D13").Value = "Au clic sur le bouton ""commander"", vérifie la validité des données ""contact"" du formulaire,  envoie au serveur les données ""contact"" et ""products"" à l'aide d'une requête de type post"
$ws.Range("E13").Value = "Renseigner un formulaire erronné et essayer de cliquer sur le bouton commander pour vérifier que l'envoi est bloqué puis renseigner correctement le formulaire et vérifier dans l'onglet Network de l'outil devtool qu'une réponse order a bien été transmise par le serveur"
$ws.Range("D14").Value = "Récupère l'order ID  envoyé par serveur et l'intègre dans l'url de la page html ""orderconfirmation"" avant de l'ouvrir"
$ws.Range("E15").Value = "Vérifier que l'orderID de l'url correspond bien à l'oderId et vérifier que le prix du panier correspond à celui contenu dans le localestorage"
$ws.Range("D16").Value = "Supprimme le contenu du locale storage si l'utlisateur quitte ou actualise la page"

$ws.Rows.Item(14).RowHeight = 60

$ws.Range("E16").Select()